$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231444597244263
$ws.Range("B1").Value = 2.489997386932373
$ws.Range("C1").Value = 4.148993968963623
$ws.Range("D1").Value = 2.768140554428101
$ws.Range("E1").Value = 1.084386825561523
